# Add a new "Resources" sheet after "Components" and populate it with the
# resource lookup table (English name / Chinese name / resource_id).

$wb = $excel.ActiveWorkbook

$componentsSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$resources = $wb.Worksheets.Add($null, $componentsSheet)
$resources.Name = "Resources"

# Header row
$resources.Range("A1").Value = "name_en"
$resources.Range("B1").Value = "name_zh"
$resources.Range("C1").Value = "resource_id"

# Chinese names, typed in row order (B2:B12)
$resources.Range("B2").Value = "铁"
$resources.Range("B3").Value = "木头"
$resources.Range("B4").Value = "皮革"
$resources.Range("B5").Value = "药草"
$resources.Range("B6").Value = "钢"
$resources.Range("B7").Value = "硬木"
$resources.Range("B8").Value = "布料"
$resources.Range("B9").Value = "油"
$resources.Range("B10").Value = "珠宝"
$resources.Range("B11").Value = "以太"
$resources.Range("B12").Value = "精华"

# resource_id values, entered alphabetically from the source lookup list
$resources.Range("C12").Value = "essence"
$resources.Range("C8").Value = "fabric"
$resources.Range("C10").Value = "gems"
$resources.Range("C5").Value = "herbs"
$resources.Range("C2").Value = "iron"
$resources.Range("C7").Value = "ironwood"
$resources.Range("C4").Value = "leather"
$resources.Range("C11").Value = "mana"
$resources.Range("C9").Value = "oils"
$resources.Range("C6").Value = "steel"
$resources.Range("C3").Value = "wood"

# English names, typed in row order (A2:A12)
$resources.Range("A2").Value = "Iron"
$resources.Range("A3").Value = "Wood"
$resources.Range("A4").Value = "Leather"
$resources.Range("A5").Value = "Herbs"
$resources.Range("A6").Value = "Steel"
$resources.Range("A7").Value = "Ironwood"
$resources.Range("A8").Value = "Fabric"
$resources.Range("A9").Value = "Oil"
$resources.Range("A10").Value = "Jewels"
$resources.Range("A11").Value = "Ether"
$resources.Range("A12").Value = "Essence"

# Make the new sheet the active tab/selection, matching the saved view state
$resources.Range("E16").Select()
$resources.Activate()
